# Apply odds updates to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 3.6
$ws.Range("G2").Value = 4.8
$ws.Range("H2").Value = 2.06
$ws.Range("J2").Value = 2.88
$ws.Range("K2").Value = 3.75
$ws.Range("N2").Value = 2.28
$ws.Range("O2").Value = 1.08
$ws.Range("P2").Value = 1.63
$ws.Range("Q2").Value = 1.94
$ws.Range("R2").Value = 1.18
$ws.Range("S2").Value = 2.68
$ws.Range("T2").Value = 1.04
$ws.Range("U2").Value = 1.04
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# Row 3
$ws.Range("P3").Value = 1.8

# Row 4
$ws.Range("F4").Value = 2.4
$ws.Range("G4").Value = 3.45
$ws.Range("H4").Value = 2.58
$ws.Range("I4").Value = 3.8
$ws.Range("J4").Value = 2.24
$ws.Range("K4").Value = 5.1
$ws.Range("L4").Value = 1.37
$ws.Range("N4").Value = 1.78
$ws.Range("P4").Value = 1.58
$ws.Range("S4").Value = 2.74
$ws.Range("V4").Value = 1.36
$ws.Range("W4").Value = 1.4

# Row 5
$ws.Range("F5").Value = 1.52
$ws.Range("G5").Value = 1.67
$ws.Range("I5").Value = 8.199999999999999
$ws.Range("L5").Value = 1.36
$ws.Range("O5").Value = 1.28
$ws.Range("P5").Value = 1.93
$ws.Range("Q5").Value = 1.82
$ws.Range("T5").Value = 1.91
$ws.Range("U5").Value = 1.83
$ws.Range("V5").Value = 1.14
$ws.Range("W5").Value = 2.48
$ws.Range("AF5").Value = 9.800000000000001

# Row 6
$ws.Range("U6").Value = 2.12
$ws.Range("V6").Value = 1.67

# Row 7
$ws.Range("G7").Value = 1.48
$ws.Range("U7").Value = 2

# Row 8
$ws.Range("AB8").Value = 970

# Row 9
$ws.Range("G9").Value = 2.16
$ws.Range("H9").Value = 3.55
$ws.Range("I9").Value = 4.6
$ws.Range("J9").Value = 3.25
$ws.Range("K9").Value = 4.2
$ws.Range("P9").Value = 1.92
$ws.Range("R9").Value = 1.38
$ws.Range("T9").Value = 1.5
$ws.Range("V9").Value = 1.28
$ws.Range("W9").Value = 1.86
